$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("E2").Value = "2016-03-14 03:45:53"
$wsZhCn.Range("H2").Value = "2016-03-14 03:46:11"

$wsDeDe.Range("E2").Value = "2016-03-14 03:45:56"
$wsDeDe.Range("H2").Value = "2016-03-14 03:46:16"
